$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("DBS")

# Update the filter-condition strings: the "ReviewType = " prefix was
# changed to "ReviewType ^i " in all three variants.
$ws.Range("B2").Value = "ReviewType ^i ,AND DataDt >= ,AND DataDt <= ,AND ProcessType %"
$ws.Range("B3").Value = "ReviewType ^i ,AND DataDt >= ,AND DataDt <= ,AND ProcessType % ,AND ProcessCount ="
$ws.Range("B4").Value = "ReviewType ^i ,AND DataDt >= ,AND DataDt <= ,AND ProcessType % ,AND ProcessCount >"

# Make "DBS" the active sheet (moves tabSelected from DBD to DBS) and move
# the active-cell selection in the frozen pane down from B11 to B12.
$ws.Activate()
$ws.Range("B12").Select()
